$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update "Curr. Est" (column C) values for the effort-estimation rows.
$ws.Range("C4").Value = 1.5
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 2
$ws.Range("C7").Value = 2

# Row 8 ("Server Socket bauen (lauffähig am Nao machen)") moved from
# "Remain" into "Effort" (fully done) and its estimate doubled.
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 0

# Row 9 ("Vom Pc aus zum Server File schicken") estimate doubled and
# also moved fully into "Effort".
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 6
$ws.Range("E9").Value = 0

$ws.Range("C10").Value = 1

# Rows 8 and 9 are now marked "Done" instead of "Open", and the
# "Open" shared string is reassigned to the new set of people still
# working the remaining tasks.
$ws.Cells.Replace("Open", "Viki, Sabina")
$ws.Range("F8").Value = "Done"
$ws.Range("F9").Value = "Done"
$ws.Range("G8").Value = "Viki, Sabina"
$ws.Range("G9").Value = "Viki, Sabina"

# Move the active selection as recorded by the author's last edit.
$ws.Range("F18").Select()

$wb.Save()
